# ============================================================================
# Applies the "final commit for demo" change:
#  - Adds two new worksheets: ProfilePage, CreateRFA
#  - Updates ErrorMessages (sheet2): replaces the now-unused "already
#    registered" message with a new "Company Name Already Exists" message,
#    and appends two new error rows (CreateRFPError / CreateRFPSupplierError)
#    plus a trailing blank-ish row.
#  - Updates SignUpPage (sheet3): adds a third "Ajay/Singh" user column (C)
#    mirroring the existing Ankit/Agarwal (B) column, including a mailto
#    hyperlink on C10.
# The order of the Value= writes below is deliberate: it reproduces the
# exact shared-string table ordering of the target workbook (new distinct
# strings get appended to the shared string table in first-use order).
# ============================================================================

$wb = $excel.ActiveWorkbook

$wsLogin  = $wb.Worksheets.Item("LogInPage")
$wsErrors = $wb.Worksheets.Item("ErrorMessages")
$wsSignUp = $wb.Worksheets.Item("SignUpPage")

# ----------------------------------------------------------------------
# 1. Add the two new sheets at the end, in order: ProfilePage, CreateRFA
# ----------------------------------------------------------------------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$wsProfile = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $lastSheet)
$wsProfile.Name = "ProfilePage"

$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$wsRFA = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $lastSheet)
$wsRFA.Name = "CreateRFA"

# ----------------------------------------------------------------------
# 2. ProfilePage!B1 -> "MaxLengthValue" (new shared string #45)
# ----------------------------------------------------------------------
$wsProfile.Range("A1").Value = "objectID"
$wsProfile.Range("B1").Value = "MaxLengthValue"
$wsLogin.Range("A1").Copy()
$wsProfile.Range("A1").PasteSpecial(-4122)
$wsLogin.Range("A1").Copy()
$wsProfile.Range("B1").PasteSpecial(-4122)

$wsProfile.Range("A2").Value = "GstField"
$wsProfile.Range("B2").Value = 15
$wsSignUp.Range("B2").Copy()
$wsProfile.Range("B2").PasteSpecial(-4122)

$wsProfile.Columns("B").ColumnWidth = 17.7109375

# ----------------------------------------------------------------------
# 3. ErrorMessages: A10/B10 -> CreateRFPError / Please Enter product name.
#    (new shared strings #46, #47)
# ----------------------------------------------------------------------
$wsErrors.Range("A10").Value = "CreateRFPError"
$wsErrors.Range("B10").Value = "Please Enter product name."

# ----------------------------------------------------------------------
# 4. CreateRFA sheet data (new shared strings #48 .. #57, in the exact
#    first-use order of the target file).
# ----------------------------------------------------------------------
$wsRFA.Range("A1").Value = "objectID"
$wsRFA.Range("A2").Value = "Mineral Water"
$wsRFA.Range("B1").Value = "Category"
$wsRFA.Range("C1").Value = "SubCategory"
$wsRFA.Range("C2").Value = "Drinking Water,Package,Packaging Machine,Treatment Equipments"
$wsRFA.Range("B3").Value = "Travel & Hotels "
$wsRFA.Range("C3").Value = "Guest Houses,OYO Rooms ,Travel & Lodging"
$wsRFA.Range("A3").Value = "Travel & Hotels"
$wsRFA.Range("D1").Value = "SuppliersName"
$wsRFA.Range("D3").Value = "ARCHER TOURS PVT LTD"
$wsRFA.Range("D2").Value = "SHIVAM ENTERPRISES"
$wsRFA.Range("B2").Value = "Mineral Water"

# Styles for CreateRFA: header row bold/filled (style 1), reuse from LogInPage!A1
$wsLogin.Range("A1").Copy()
$wsRFA.Range("A1:D1").PasteSpecial(-4122)

# Column widths on CreateRFA
$wsRFA.Columns("A").ColumnWidth = 21.5546875
$wsRFA.Columns("B").ColumnWidth = 19.44140625
$wsRFA.Columns("C").ColumnWidth = 53.7109375
$wsRFA.Columns("D").ColumnWidth = 20.109375

# ----------------------------------------------------------------------
# 5. ErrorMessages: A11/B11 -> CreateRFPSupplierError / Please select
#    atleast one supplier. (new shared strings #58, #59)
# ----------------------------------------------------------------------
$wsErrors.Range("A11").Value = "CreateRFPSupplierError"
$wsErrors.Range("B11").Value = "Please select atleast one supplier."

# ----------------------------------------------------------------------
# 6. SignUpPage column C (Ajay / Singh user), new shared strings #60, #61
#    for the name, remaining values reuse existing shared strings.
# ----------------------------------------------------------------------
$wsSignUp.Range("C2").Value = "Ajay"
$wsSignUp.Range("C3").Value = "Singh"
$wsSignUp.Range("C4").Value = "Sam Technology"
$wsSignUp.Range("C5").Value = "Noida sector 62"
$wsSignUp.Range("C6").Value = "Uttar Pradesh"
$wsSignUp.Range("C7").Value = "Noida"
$wsSignUp.Range("C8").Value = 201301
$wsSignUp.Range("C9").Value = "01234AFD12"
$wsSignUp.Range("C10").Value = "hprankit@gmail.com"
$wsSignUp.Range("C11").Value = 9568989975
$wsSignUp.Range("C12").Value = "qwerty11"

# Reuse formatting from column B for column C
$wsSignUp.Range("B2:B9").Copy()
$wsSignUp.Range("C2").PasteSpecial(-4122)
$wsSignUp.Range("B10").Copy()
$wsSignUp.Range("C10").PasteSpecial(-4122)
$wsSignUp.Range("B11:B12").Copy()
$wsSignUp.Range("C11").PasteSpecial(-4122)

# Hyperlink for C10 (mailto), added after the style paste so the hyperlink
# style id matches the existing B10 style id (4) instead of minting a new one.
$wsSignUp.Hyperlinks.Add($wsSignUp.Range("C10"), "mailto:hprankit@gmail.com")
$wsSignUp.Range("B10").Copy()
$wsSignUp.Range("C10").PasteSpecial(-4122)

# ----------------------------------------------------------------------
# 7. ErrorMessages: A12 -> " " (reuses existing shared string #4) and
#    B8 -> "Company Name Already Exists..." (new shared string #62,
#    replaces the now-orphaned "already exists" message so it drops out
#    of the shared-string table on save).
# ----------------------------------------------------------------------
$wsErrors.Range("A12").Value = " "
$wsErrors.Range("B8").Value = "Company Name Already Exists. Please try With Different Name."

# ----------------------------------------------------------------------
# 8. Selections / active cells, matching the target sheetViews. Set the
#    non-active sheets' selections first, then re-activate ErrorMessages
#    (the originally active tab) last so tabSelected stays correct.
# ----------------------------------------------------------------------
$wsProfile.Range("A1:B1").Select()
$wsRFA.Range("C13").Select()
$wsSignUp.Range("G15").Select()
$wsErrors.Activate()
$wsErrors.Range("B9").Select()

Write-Host "edit complete"
